# Updated cryptos list on Thu Aug 15 09:33:15 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row in the cryptocurrency tracking sheet.  Values that look like
# plain numbers (e.g. "520.11") are forced to Text format before being
# written so Excel keeps them as literal strings (matching the source
# data's formatting, e.g. trailing zeros / exact decimals) instead of
# re-interpreting them as numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.360.93"
$ws.Range("E2").Value = "  -3.76%  "
$ws.Range("D3").Value = "2.616.91"
$ws.Range("E3").Value = "  -3.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.11"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.97"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.61"
$ws.Range("E9").Value = "  -3.48%  "
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").Value = "3.077.26"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").Value = "58.321.32"
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.97"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "2.611.26"
$ws.Range("E17").Value = "  -10.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "336.45"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("E19").Value = "  -2.88%  "
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.49"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.415"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").Value = "  -2.24%  "
$ws.Range("E28").Value = "  -3.86%  "
$ws.Range("E29").Value = "  -3.56%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.78"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.98"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("E34").Value = "  -3.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.18"
$ws.Range("E35").Value = "  -3.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.884"
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.855"
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.34"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("E39").Value = "  -6.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.603"
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "269.09"
$ws.Range("E44").Value = "  -4.75%  "
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.14"
$ws.Range("E46").Value = "  -5.18%  "
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("D48").Value = "2.034.42"
$ws.Range("E48").Value = "  -4.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0228"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("E50").Value = "  -5.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.29"
$ws.Range("E51").Value = "  -5.09%  "
